$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B to make room for "NumeroFactura".
# This shifts old B..G (EPS, Valor, Vigencia, Estado, Mes, Observaciones) to C..H.
$ws.Columns.Item(2).Insert()

# Fill in the new "NumeroFactura" column.
$ws.Range("B1").Value = "NumeroFactura"
$ws.Range("B2").Value = "FAC001"
$ws.Range("B3").Value = "FAC002"
$ws.Range("B4").Value = "FAC003"
$ws.Range("B5").Value = "FAC004"

# After the insert, column order is: ID, (new), EPS, Valor, Vigencia, Estado, Mes, Observaciones
# Target order is:                   ID, NumeroFactura, Valor, EPS, Vigencia, Estado, Mes, Observaciones
# So swap columns C (EPS) and D (Valor).
$ws.Range("C1").Value = "Valor"
$ws.Range("D1").Value = "EPS"

$ws.Range("C2").Value = 150000
$ws.Range("D2").Value = "Sanitas"

$ws.Range("C3").Value = 250000
$ws.Range("D3").Value = "Compensar"

$ws.Range("C4").Value = 300000
$ws.Range("D4").Value = "Sura"

$ws.Range("C5").Value = 100000
$ws.Range("D5").Value = "Nueva EPS"

# The "Estado" column (now F) shifts down one row (rotates).
$ws.Range("F2").Value = "Pendiente"
$ws.Range("F3").Value = "Auditada"
$ws.Range("F4").Value = "Subsanada"
$ws.Range("F5").Value = "Radicada"

Write-Host "Workbook updated"
